# Append the new profit-allocation row for 2025-09-25 (row 24) to the sheet.
#
# Column A holds a literal date string like "09/25/2025" (stored as text in
# the source file, not as a real Excel date). Simply assigning the string to
# .Value would make Excel auto-detect it as a date and store a date serial
# number instead, so we briefly mark the cell as Text ("@") before writing
# the value, then clear the formatting back to the sheet's default so the
# new row matches the look of the existing rows above it.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 24

$dateCell = $ws.Range("A$row")
$dateCell.NumberFormat = "@"
$dateCell.Value = "09/25/2025"
$dateCell.ClearFormats()

$ws.Range("B$row").Value = 0.1344413231227061
$ws.Range("C$row").Value = 0.8655586768772939
